$wb = $excel.ActiveWorkbook

# 1. Reorder worksheet tabs: "review_info" moves in front of "hotel_info"
#    (workbook originally had hotel_info, review_info -> now review_info, hotel_info).
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($hotelSheet)

# 2. On hotel_info, insert a new "State" column between "Hotel_Name" (B) and
#    "City" (C), shifting City/Zip/... one column to the right, and fill in
#    the value for the existing data row.
$ws = $wb.Worksheets.Item("hotel_info")
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("C1").Value = "State"
$ws.Range("C2").Value = "Louisiana"
